$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.034.49"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "3.509.98"
$ws.Range("E3").Value = "  -2.14%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'573.20"
$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("D6").Value = "'184.62"
$ws.Range("E6").Value = "  -3.61%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.613"
$ws.Range("E7").Value = "  -3.21%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.491.20"
$ws.Range("E8").Value = "  -2.51%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("D11").Value = "'0.650"
$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").Value = "'54.23"
$ws.Range("E12").Value = "  -2.99%  "

$ws.Range("D13").Value = "'0.0000301"
$ws.Range("E13").Value = "  -1.63%  "

$ws.Range("D14").Value = "'9.44"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").Value = "4.069.25"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").Value = "'19.37"
$ws.Range("E16").Value = "  -2.87%  "

$ws.Range("D17").Value = "3.509.21"
$ws.Range("E17").Value = "  -2.19%  "

$ws.Range("D18").Value = "68.943.49"
$ws.Range("E18").Value = "  -1.47%  "

$ws.Range("D19").Value = "'12.30"
$ws.Range("E19").Value = "  -3.13%  "

$ws.Range("E20").Value = "  -1.29%  "

$ws.Range("D21").Value = "'542.15"
$ws.Range("E21").Value = "  +13.91%  "

$ws.Range("E22").Value = "  -3.34%  "

$ws.Range("D23").Value = "'18.30"
$ws.Range("E23").Value = "  -5.45%  "

$ws.Range("D24").Value = "'4.96"
$ws.Range("E24").Value = "  -1.45%  "

$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("D26").Value = "'94.05"
$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.92"
$ws.Range("E27").Value = "  -2.73%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.92"
$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("E29").Value = "  -3.45%  "

$ws.Range("D30").Value = "'31.68"
$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("D31").Value = "'7.24"
$ws.Range("E31").Value = "  -5.65%  "

$ws.Range("D32").Value = "'12.70"
$ws.Range("E32").Value = "  +3.71%  "

$ws.Range("D33").Value = "'64.76"
$ws.Range("E33").Value = "  -2.59%  "

$ws.Range("E34").Value = "  -4.46%  "

$ws.Range("D35").Value = "'551.96"
$ws.Range("E35").Value = "  -6.34%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'38.01"
$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'3.08"
$ws.Range("E37").Value = "  +7.80%  "

$ws.Range("E38").Value = "  +1.16%  "

$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("E40").Value = "  -4.88%  "

$ws.Range("D41").Value = "'3.36"
$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("D42").Value = "'3.08"
$ws.Range("E42").Value = "  -4.85%  "

$ws.Range("E43").Value = "  -3.55%  "

$ws.Range("D44").Value = "3.291.70"
$ws.Range("E44").Value = "  +1.74%  "

$ws.Range("D45").Value = "'2.99"
$ws.Range("E45").Value = "  -3.57%  "

$ws.Range("D46").Value = "'0.0444"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'3.46"
$ws.Range("E47").Value = "  +3.35%  "

$ws.Range("E48").Value = "  -2.85%  "

$ws.Range("D49").Value = "'8.89"
$ws.Range("E49").Value = "  -5.98%  "

$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("D51").Value = "'137.90"
$ws.Range("E51").Value = "  +2.75%  "
